# Auto-generated cell updates applying the Ixion_Profits.xlsx diff
$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("N26").ClearContents()
$ws.Range("L26").Value = 0
$ws.Range("H26").Value = 0
$ws.Range("J26").Value = 0
$ws.Range("M34").Value = -1943.8333
$ws.Range("I34").Value = 2146.8333
$ws.Range("K34").Value = 2146.8333
$ws.Range("H34").Value = 3854.4285
$ws.Range("I36").Value = 2146.8333
$ws.Range("K36").Value = 2146.8333
$ws.Range("M36").Value = -1431.8333
$ws.Range("H36").Value = 3854.4285
$ws.Range("H62").Value = 1955.6111
$ws.Range("J62").Value = 2624.5
$ws.Range("N62").Value = -3872.5
$ws.Range("K62").Value = 1764.5
$ws.Range("I62").Value = 1764.5
$ws.Range("M62").Value = -1140.5
$ws.Range("L62").Value = 2624.5
$ws.Range("H64").Value = 3697.963
$ws.Range("K64").Value = 3889.7437
$ws.Range("J64").Value = 3199.3333
$ws.Range("N64").Value = -3695.3333
$ws.Range("L64").Value = 3199.3333
$ws.Range("M64").Value = -3641.7437
$ws.Range("I64").Value = 3889.7437
$ws.Range("L65").Value = 13122.5
$ws.Range("N65").Value = -19362.5
$ws.Range("M65").Value = -5702.5
$ws.Range("K65").Value = 8822.5
$ws.Range("I65").Value = 1764.5
$ws.Range("J65").Value = 2624.5
$ws.Range("H65").Value = 1955.6111
$ws.Range("K67").Value = 3889.7437
$ws.Range("M67").Value = -3031.7437
$ws.Range("N67").Value = -4915.3333
$ws.Range("H67").Value = 3697.963
$ws.Range("J67").Value = 3199.3333
$ws.Range("L67").Value = 3199.3333
$ws.Range("I67").Value = 3889.7437
$ws.Range("K74").Value = 3914.2856
$ws.Range("M74").Value = -2978.2856
$ws.Range("I74").Value = 3914.2856
$ws.Range("H74").Value = 3800
$ws.Range("H77").Value = 3800
$ws.Range("I77").Value = 3914.2856
$ws.Range("M77").Value = -14891.428
$ws.Range("K77").Value = 19571.428
$ws.Range("K80").Value = 1117.00002
$ws.Range("M80").Value = -119.0000199999999
$ws.Range("L80").Value = 429
$ws.Range("H80").Value = 315
$ws.Range("J80").Value = 143
$ws.Range("I80").Value = 372.33334
$ws.Range("N80").Value = -2425
$ws.Range("K83").Value = 3351.00006
$ws.Range("N83").Value = -11271
$ws.Range("L83").Value = 1287
$ws.Range("I83").Value = 372.33334
$ws.Range("M83").Value = 1640.99994
$ws.Range("H83").Value = 315
$ws.Range("J83").Value = 143
$ws.Range("H129").Value = 1037.2285
$ws.Range("N129").Value = -13273.1383
$ws.Range("I129").Value = 337.6
$ws.Range("M129").Value = 3987.2
$ws.Range("L129").Value = 3273.1383
$ws.Range("K129").Value = 1012.8
$ws.Range("J129").Value = 1091.0461
$ws.Range("H133").Value = 39945
$ws.Range("J133").Value = 39945
$ws.Range("N133").Value = -50065
$ws.Range("L133").Value = 39945
$ws.Range("K137").Value = 4562.4288
$ws.Range("M137").Value = -2012.4288
$ws.Range("H137").Value = 1504.2307
$ws.Range("I137").Value = 1520.8096
$ws.Range("K138").Value = 2584.9617
$ws.Range("H138").Value = 1820.3651
$ws.Range("N138").Value = -17762.162
$ws.Range("L138").Value = 7482.162
$ws.Range("I138").Value = 861.6539
$ws.Range("M138").Value = 2555.0383
$ws.Range("J138").Value = 2494.054

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("N61").Value = -480230.62
$ws.Range("M61").Value = -1666.5385
$ws.Range("K61").Value = 1878.5385
$ws.Range("L61").Value = 479806.62
$ws.Range("J61").Value = 479806.62
$ws.Range("I61").Value = 1878.5385
$ws.Range("H61").Value = 215420.88
$ws.Range("H136").Value = 215420.88
$ws.Range("M136").Value = -3085.6155
$ws.Range("J136").Value = 479806.62
$ws.Range("K136").Value = 5635.6155
$ws.Range("N136").Value = -1444519.86
$ws.Range("I136").Value = 1878.5385
$ws.Range("L136").Value = 1439419.86

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("M20").Value = -1116.2307
$ws.Range("K20").Value = 1363.2307
$ws.Range("I20").Value = 1363.2307
$ws.Range("N20").Value = -44269
$ws.Range("H20").Value = 14756.421
$ws.Range("L20").Value = 43775
$ws.Range("J20").Value = 43775
$ws.Range("M94").Value = 172
$ws.Range("H94").Value = 1186
$ws.Range("I94").Value = 279
$ws.Range("K94").Value = 279

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("I31").Value = 2050.4
$ws.Range("K31").Value = 2050.4
$ws.Range("L31").Value = 4130.2
$ws.Range("H31").Value = 3090.3
$ws.Range("N31").Value = -4720.2
$ws.Range("M31").Value = -1755.4
$ws.Range("J31").Value = 4130.2
$ws.Range("M34").Value = -1848.4
$ws.Range("I34").Value = 2050.4
$ws.Range("L34").Value = 4130.2
$ws.Range("K34").Value = 2050.4
$ws.Range("H34").Value = 3090.3
$ws.Range("N34").Value = -4534.2
$ws.Range("J34").Value = 4130.2
$ws.Range("K132").Value = 6157.875
$ws.Range("M132").Value = -3627.875
$ws.Range("I132").Value = 2052.625
$ws.Range("H132").Value = 2438.32
$ws.Range("M134").Value = -6524.4828
$ws.Range("J134").Value = 1857
$ws.Range("H134").Value = 2944.8064
$ws.Range("L134").Value = 5571
$ws.Range("K134").Value = 9059.4828
$ws.Range("N134").Value = -10641
$ws.Range("I134").Value = 3019.8276

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H69").Value = 1103.875
$ws.Range("M69").Value = -690.8000000000002
$ws.Range("I69").Value = 500.6
$ws.Range("N69").Value = -7949.999899999999
$ws.Range("J69").Value = 2109.3333
$ws.Range("L69").Value = 6327.999899999999
$ws.Range("K69").Value = 1501.8
$ws.Range("J72").Value = 2109.3333
$ws.Range("I72").Value = 500.6
$ws.Range("N72").Value = -27095.9997
$ws.Range("M72").Value = -449.4000000000005
$ws.Range("L72").Value = 18983.9997
$ws.Range("H72").Value = 1103.875
$ws.Range("K72").Value = 4505.400000000001
$ws.Range("N74").Value = -62122
$ws.Range("L74").Value = 60000
$ws.Range("H74").Value = 11000
$ws.Range("J74").Value = 20000
$ws.Range("N75").Value = -214288216
$ws.Range("L75").Value = 214286220
$ws.Range("H75").Value = 71428740
$ws.Range("J75").Value = 71428740
$ws.Range("H77").Value = 11000
$ws.Range("J77").Value = 20000
$ws.Range("N77").Value = -190608
$ws.Range("L77").Value = 180000
$ws.Range("L78").Value = 642858660
$ws.Range("H78").Value = 71428740
$ws.Range("N78").Value = -642868644
$ws.Range("J78").Value = 71428740
$ws.Range("M113").Value = -7499459.600000001
$ws.Range("L113").Value = 1876595.25
$ws.Range("H113").Value = 1667204.8
$ws.Range("N113").Value = -1880935.25
$ws.Range("J113").Value = 625531.75
$ws.Range("K113").Value = 7501629.600000001
$ws.Range("I113").Value = 2500543.2

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 333553.97
$ws.Range("L102").Value = 1733.5238
$ws.Range("N102").Value = -4977.5238
$ws.Range("K102").Value = 565828.25
$ws.Range("I102").Value = 565828.25
$ws.Range("M102").Value = -564206.25
$ws.Range("J102").Value = 1733.5238
$ws.Range("M113").Value = -76921950
$ws.Range("L113").Value = 2880
$ws.Range("H113").Value = 55557110
$ws.Range("N113").Value = -7220
$ws.Range("J113").Value = 2880
$ws.Range("K113").Value = 76924120
$ws.Range("I113").Value = 76924120
$ws.Range("J134").Value = 0
$ws.Range("H134").Value = 0
$ws.Range("L134").Value = 0
$ws.Range("N134").ClearContents()
$ws.Range("H135").Value = 32756
$ws.Range("J135").Value = 32756
$ws.Range("N135").Value = -42896
$ws.Range("L135").Value = 32756

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("I93").Value = 3741
$ws.Range("H93").Value = 2980.75
$ws.Range("K93").Value = 3741
$ws.Range("J93").Value = 700
$ws.Range("L93").Value = 700
$ws.Range("M93").Value = -2493
$ws.Range("N93").Value = -3196

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("L43").Value = 9666.666999999999
$ws.Range("H43").Value = 9666.666999999999
$ws.Range("J43").Value = 9666.666999999999
$ws.Range("N43").Value = -9964.666999999999
$ws.Range("K107").Value = 166668390
$ws.Range("J107").Value = 1201.2
$ws.Range("L107").Value = 3603.6
$ws.Range("I107").Value = 55556130
$ws.Range("N107").Value = -7443.6
$ws.Range("H107").Value = 43478972
$ws.Range("M107").Value = -166666470
$ws.Range("L126").Value = 3975
$ws.Range("H126").Value = 788
$ws.Range("K126").Value = 2041.8
$ws.Range("J126").Value = 1325
$ws.Range("N126").Value = -8915
$ws.Range("M126").Value = 428.1999999999998
$ws.Range("I126").Value = 680.6
$ws.Range("H136").Value = 2874.3225
$ws.Range("M136").Value = -8049.799800000001
$ws.Range("J136").Value = 2256.5625
$ws.Range("K136").Value = 10599.7998
$ws.Range("I136").Value = 3533.2666
$ws.Range("L136").Value = 6769.6875

Write-Output "Applied all cell updates: 222 sets, 2 clears"